$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.924095199937824
$ws.Cells.Item(2, 3).Value = 0.2797067832110969
$ws.Cells.Item(2, 4).Value = 0.07882263183346083
$ws.Cells.Item(2, 5).Value = 0.126489902436834
$ws.Cells.Item(2, 7).Value = 0.5774085957165198
$ws.Cells.Item(2, 8).Value = 0.6919621358620986
$ws.Cells.Item(2, 12).Value = 0.1882037178545133
$ws.Cells.Item(2, 13).Value = 0.2088288491555161
$ws.Cells.Item(2, 14).Value = 1.319205680391462
$ws.Cells.Item(2, 15).Value = 2.516099098855335

$ws.Cells.Item(3, 2).Value = 0.8410414654981082
$ws.Cells.Item(3, 3).Value = 0.2705063221765016
$ws.Cells.Item(3, 4).Value = 0.07146014508909104
$ws.Cells.Item(3, 5).Value = 0.1276255021204991
$ws.Cells.Item(3, 7).Value = 0.5736217929417933
$ws.Cells.Item(3, 8).Value = 0.6944288376109995
$ws.Cells.Item(3, 12).Value = 0.1855205434652305
$ws.Cells.Item(3, 13).Value = 0.1956210123980142
$ws.Cells.Item(3, 14).Value = 1.331747790706743
$ws.Cells.Item(3, 15).Value = 2.512923386086953

$ws.Cells.Item(4, 2).Value = 0.7901881031667415
$ws.Cells.Item(4, 3).Value = 0.2648167793937546
$ws.Cells.Item(4, 4).Value = 0.06697365006199618
$ws.Cells.Item(4, 5).Value = 0.128365758579517
$ws.Cells.Item(4, 7).Value = 0.5717304701314561
$ws.Cells.Item(4, 8).Value = 0.6962938461888086
$ws.Cells.Item(4, 12).Value = 0.1839638322917736
$ws.Cells.Item(4, 13).Value = 0.187576668066221
$ws.Cells.Item(4, 14).Value = 1.339954314033129
$ws.Cells.Item(4, 15).Value = 2.512615449442706

$ws.Cells.Item(5, 2).Value = 0.7695018796119655
$ws.Cells.Item(5, 3).Value = 0.2624882315663797
$ws.Cells.Item(5, 4).Value = 0.06515394383714579
$ws.Cells.Item(5, 5).Value = 0.128678241657938
$ws.Cells.Item(5, 7).Value = 0.5710686682970163
$ws.Cells.Item(5, 8).Value = 0.6971419998727555
$ws.Cells.Item(5, 12).Value = 0.1833523414498828
$ws.Cells.Item(5, 13).Value = 0.1843151571724704
$ws.Cells.Item(5, 14).Value = 1.3434257426565
$ws.Cells.Item(5, 15).Value = 2.512902570839856

$ws.Cells.Item(6, 2).Value = 0.7660692151328874
$ws.Cells.Item(6, 3).Value = 0.2621009776851224
$ws.Cells.Item(6, 4).Value = 0.06485230108532392
$ws.Cells.Item(6, 5).Value = 0.1287307832281361
$ws.Cells.Item(6, 7).Value = 0.5709653511059685
$ws.Cells.Item(6, 8).Value = 0.6972881598055238
$ws.Cells.Item(6, 12).Value = 0.1832521875360982
$ws.Cells.Item(6, 13).Value = 0.183774595740843
$ws.Cells.Item(6, 14).Value = 1.34400985543013
$ws.Cells.Item(6, 15).Value = 2.512975159563041

$ws.Cells.Item(7, 2).Value = 0.7899089703362847
$ws.Cells.Item(7, 3).Value = 0.2647854161249796
$ws.Cells.Item(7, 4).Value = 0.06694907413731244
$ws.Cells.Item(7, 5).Value = 0.1283699290001072
$ws.Cells.Item(7, 7).Value = 0.5717211039996641
$ws.Cells.Item(7, 8).Value = 0.6963049277470503
$ws.Cells.Item(7, 12).Value = 0.1839554927895648
$ws.Cells.Item(7, 13).Value = 0.1875326146124152
$ws.Cells.Item(7, 14).Value = 1.340000615806893
$ws.Cells.Item(7, 15).Value = 2.512617651399324

$ws.Cells.Item(8, 2).Value = 0.8954295218493371
$ws.Cells.Item(8, 3).Value = 0.276542942164582
$ws.Cells.Item(8, 4).Value = 0.07627696748198787
$ws.Cells.Item(8, 5).Value = 0.1268725408429585
$ws.Cells.Item(8, 7).Value = 0.5760127786172689
$ws.Cells.Item(8, 8).Value = 0.6927399342349219
$ws.Cells.Item(8, 12).Value = 0.1872597531980702
$ws.Cells.Item(8, 13).Value = 0.2042613447965991
$ws.Cells.Item(8, 14).Value = 1.323425296449237
$ws.Cells.Item(8, 15).Value = 2.514663148161333

$ws.Cells.Item(9, 2).Value = 1.103437601374935
$ws.Cells.Item(9, 3).Value = 0.299273044204682
$ws.Cells.Item(9, 4).Value = 0.09484041228483875
$ws.Cells.Item(9, 5).Value = 0.1242767632346302
$ws.Cells.Item(9, 7).Value = 0.5878789111850011
$ws.Cells.Item(9, 8).Value = 0.6885290443847083
$ws.Cells.Item(9, 12).Value = 0.1944576435740828
$ws.Cells.Item(9, 13).Value = 0.2375776342353149
$ws.Cells.Item(9, 14).Value = 1.294930137318012
$ws.Cells.Item(9, 15).Value = 2.531719269780979

$ws.Cells.Item(10, 2).Value = 1.256877612054382
$ws.Cells.Item(10, 3).Value = 0.3157678197280802
$ws.Cells.Item(10, 4).Value = 0.1086471137551541
$ws.Cells.Item(10, 5).Value = 0.122576503346794
$ws.Cells.Item(10, 7).Value = 0.5987136755283018
$ws.Cells.Item(10, 8).Value = 0.6871301786048463
$ws.Cells.Item(10, 12).Value = 0.2001819588782325
$ws.Cells.Item(10, 13).Value = 0.262360020095997
$ws.Cells.Item(10, 14).Value = 1.27643472567506
$ws.Cells.Item(10, 15).Value = 2.55223199841933

$ws.Cells.Item(11, 2).Value = 1.326807028301175
$ws.Cells.Item(11, 3).Value = 0.323226007018377
$ws.Cells.Item(11, 4).Value = 0.1149652884714527
$ws.Cells.Item(11, 5).Value = 0.1218477621631768
$ws.Cells.Item(11, 7).Value = 0.6041053068771447
$ws.Cells.Item(11, 8).Value = 0.686861909193496
$ws.Cells.Item(11, 12).Value = 0.2028804385164307
$ws.Cells.Item(11, 13).Value = 0.2736990962176478
$ws.Cells.Item(11, 14).Value = 1.268549729037147
$ws.Cells.Item(11, 15).Value = 2.563303521976252

$ws.Cells.Item(12, 2).Value = 1.353304952415044
$ws.Cells.Item(12, 3).Value = 0.3260435700693165
$ws.Cells.Item(12, 4).Value = 0.1173632179967399
$ws.Cells.Item(12, 5).Value = 0.1215782252691779
$ws.Cells.Item(12, 7).Value = 0.6062137332341138
$ws.Cells.Item(12, 8).Value = 0.6868132465388754
$ws.Cells.Item(12, 12).Value = 0.203915820268108
$ws.Cells.Item(12, 13).Value = 0.2780021538451223
$ws.Cells.Item(12, 14).Value = 1.265639868302046
$ws.Cells.Item(12, 15).Value = 2.567746674955259

$ws.Cells.Item(13, 2).Value = 1.347597408061858
$ws.Cells.Item(13, 3).Value = 0.3254370576690917
$ws.Cells.Item(13, 4).Value = 0.1168465421247333
$ws.Cells.Item(13, 5).Value = 0.1216359893948011
$ws.Cells.Item(13, 7).Value = 0.6057566751496495
$ws.Cells.Item(13, 8).Value = 0.6868213731234079
$ws.Cells.Item(13, 12).Value = 0.2036922316171825
$ws.Cells.Item(13, 13).Value = 0.277075007842555
$ws.Cells.Item(13, 14).Value = 1.266263178049371
$ws.Cells.Item(13, 15).Value = 2.566778610748145

$ws.Cells.Item(14, 2).Value = 1.328986692027854
$ws.Cells.Item(14, 3).Value = 0.3234579446696841
$ws.Cells.Item(14, 4).Value = 0.1151624599336571
$ws.Cells.Item(14, 5).Value = 0.1218254585509522
$ws.Cells.Item(14, 7).Value = 0.6042774299563973
$ws.Cells.Item(14, 8).Value = 0.6868568450145602
$ws.Cells.Item(14, 12).Value = 0.2029653492409267
$ws.Cells.Item(14, 13).Value = 0.2740529282973867
$ws.Cells.Item(14, 14).Value = 1.268308809571252
$ws.Cells.Item(14, 15).Value = 2.563664038183987

$ws.Cells.Item(15, 2).Value = 1.317589287404928
$ws.Cells.Item(15, 3).Value = 0.3222448045583235
$ws.Cells.Item(15, 4).Value = 0.1141316105075703
$ws.Cells.Item(15, 5).Value = 0.1219423499086592
$ws.Cells.Item(15, 7).Value = 0.6033800453270999
$ws.Cells.Item(15, 8).Value = 0.6868854648560614
$ws.Cells.Item(15, 12).Value = 0.2025218725958666
$ws.Cells.Item(15, 13).Value = 0.2722030084672156
$ws.Cells.Item(15, 14).Value = 1.269571717466242
$ws.Cells.Item(15, 15).Value = 2.561788920004432

$ws.Cells.Item(16, 2).Value = 1.252310036525046
$ws.Cells.Item(16, 3).Value = 0.3152794824583793
$ws.Cells.Item(16, 4).Value = 0.1082349587113498
$ws.Cells.Item(16, 5).Value = 0.1226250273631129
$ws.Cells.Item(16, 7).Value = 0.5983706484236251
$ws.Cells.Item(16, 8).Value = 0.6871551169939636
$ws.Cells.Item(16, 12).Value = 0.2000075027778081
$ws.Cells.Item(16, 13).Value = 0.2616202833108616
$ws.Cells.Item(16, 14).Value = 1.276960666269936
$ws.Cells.Item(16, 15).Value = 2.551543497008623

$ws.Cells.Item(17, 2).Value = 1.212295356397703
$ws.Cells.Item(17, 3).Value = 0.3109947442065675
$ws.Cells.Item(17, 4).Value = 0.1046271358196265
$ws.Cells.Item(17, 5).Value = 0.1230552745516249
$ws.Cells.Item(17, 7).Value = 0.595416216697302
$ws.Cells.Item(17, 8).Value = 0.6874148111188134
$ws.Cells.Item(17, 12).Value = 0.1984891711306744
$ws.Cells.Item(17, 13).Value = 0.2551447328972856
$ws.Cells.Item(17, 14).Value = 1.281628945010539
$ws.Cells.Item(17, 15).Value = 2.54570422892877

$ws.Cells.Item(18, 2).Value = 1.189292173093975
$ws.Cells.Item(18, 3).Value = 0.3085260156953495
$ws.Cells.Item(18, 4).Value = 0.1025555358821322
$ws.Cells.Item(18, 5).Value = 0.1233069506979394
$ws.Cells.Item(18, 7).Value = 0.5937604601929962
$ws.Cells.Item(18, 8).Value = 0.6875988261564885
$ws.Cells.Item(18, 12).Value = 0.1976247613567352
$ws.Cells.Item(18, 13).Value = 0.2514263403722197
$ws.Cells.Item(18, 14).Value = 1.284363775325737
$ws.Cells.Item(18, 15).Value = 2.542509408329551

$ws.Cells.Item(19, 2).Value = 1.181505827328238
$ws.Cells.Item(19, 3).Value = 0.3076894199986384
$ws.Cells.Item(19, 4).Value = 0.1018547333271016
$ws.Cells.Item(19, 5).Value = 0.1233928871345107
$ws.Cells.Item(19, 7).Value = 0.5932073246213037
$ws.Cells.Item(19, 8).Value = 0.6876670810867154
$ws.Cells.Item(19, 12).Value = 0.197333616427116
$ws.Cells.Item(19, 13).Value = 0.2501684233003019
$ws.Cells.Item(19, 14).Value = 1.285298288378939
$ws.Cells.Item(19, 15).Value = 2.541455813607115

$ws.Cells.Item(20, 2).Value = 1.216553734274498
$ws.Cells.Item(20, 3).Value = 0.3114513040512179
$ws.Cells.Item(20, 4).Value = 0.1050108300343311
$ws.Cells.Item(20, 5).Value = 0.1230090384144162
$ws.Cells.Item(20, 7).Value = 0.5957262119905664
$ws.Cells.Item(20, 8).Value = 0.6873835805528046
$ws.Cells.Item(20, 12).Value = 0.1986498800005876
$ws.Cells.Item(20, 13).Value = 0.2558334289915862
$ws.Cells.Item(20, 14).Value = 1.281126848851883
$ws.Cells.Item(20, 15).Value = 2.546308876498301

$ws.Cells.Item(21, 2).Value = 1.33445265514797
$ws.Cells.Item(21, 3).Value = 0.32403944103919
$ws.Cells.Item(21, 4).Value = 0.1156569699299581
$ws.Cells.Item(21, 5).Value = 0.1217696326544546
$ws.Cells.Item(21, 7).Value = 0.6047101076249675
$ws.Cells.Item(21, 8).Value = 0.6868449897228857
$ws.Cells.Item(21, 12).Value = 0.203178485576629
$ws.Cells.Item(21, 13).Value = 0.2749403384650222
$ws.Cells.Item(21, 14).Value = 1.267705894956784
$ws.Cells.Item(21, 15).Value = 2.564572059286832

$ws.Cells.Item(22, 2).Value = 1.411605882030585
$ws.Cells.Item(22, 3).Value = 0.3322274412425941
$ws.Cells.Item(22, 4).Value = 0.1226461559016059
$ws.Cells.Item(22, 5).Value = 0.120997035324188
$ws.Cells.Item(22, 7).Value = 0.6109706522229885
$ws.Cells.Item(22, 8).Value = 0.6868014757281173
$ws.Cells.Item(22, 12).Value = 0.2062170049044312
$ws.Cells.Item(22, 13).Value = 0.2874813075087488
$ws.Cells.Item(22, 14).Value = 1.259377608257928
$ws.Cells.Item(22, 15).Value = 2.577969006641183

$ws.Cells.Item(23, 2).Value = 1.370419096376622
$ws.Cells.Item(23, 3).Value = 0.3278609817402014
$ws.Cells.Item(23, 4).Value = 0.1189130335504842
$ws.Cells.Item(23, 5).Value = 0.1214059633594387
$ws.Cells.Item(23, 7).Value = 0.6075936257632009
$ws.Cells.Item(23, 8).Value = 0.6867964749266804
$ws.Cells.Item(23, 12).Value = 0.2045880973723797
$ws.Cells.Item(23, 13).Value = 0.2807831322068424
$ws.Cells.Item(23, 14).Value = 1.263782029842361
$ws.Cells.Item(23, 15).Value = 2.5706850090464

$ws.Cells.Item(24, 2).Value = 1.21462851595993
$ws.Cells.Item(24, 3).Value = 0.3112449100769652
$ws.Cells.Item(24, 4).Value = 0.1048373538266816
$ws.Cells.Item(24, 5).Value = 0.1230299283110554
$ws.Cells.Item(24, 7).Value = 0.5955859298572364
$ws.Cells.Item(24, 8).Value = 0.6873975917594635
$ws.Cells.Item(24, 12).Value = 0.1985771970426384
$ws.Cells.Item(24, 13).Value = 0.2555220554926407
$ws.Cells.Item(24, 14).Value = 1.281353687739319
$ws.Cells.Item(24, 15).Value = 2.546035009950856

$ws.Cells.Item(25, 2).Value = 1.047054510119096
$ws.Cells.Item(25, 3).Value = 0.2931595183826801
$ws.Cells.Item(25, 4).Value = 0.08978916972237982
$ws.Cells.Item(25, 5).Value = 0.1249426074706568
$ws.Cells.Item(25, 7).Value = 0.584298054283849
$ws.Cells.Item(25, 8).Value = 0.6893705705227262
$ws.Cells.Item(25, 12).Value = 0.192433720287859
$ws.Cells.Item(25, 13).Value = 0.2285106602816001
$ws.Cells.Item(25, 14).Value = 1.302210073419943
$ws.Cells.Item(25, 15).Value = 2.525705578130641
